$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.765.60'
$ws.Range('E2').Value = '  +3.77%  '
$ws.Range('D3').Value = '3.691.55'
$ws.Range('E3').Value = '  +8.15%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '590.08'
$ws.Range('E5').Value = '  +1.55%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '180.95'
$ws.Range('E6').Value = '  +1.71%  '
$ws.Range('D7').Value = '3.676.96'
$ws.Range('E7').Value = '  +8.00%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.618'
$ws.Range('E8').Value = '  +4.45%  '
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.203'
$ws.Range('E10').Value = '  +2.66%  '
$ws.Range('E11').Value = '  +5.18%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '50.09'
$ws.Range('E12').Value = '  +3.90%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000288'
$ws.Range('E13').Value = '  +2.58%  '
$ws.Range('D14').Value = '4.283.62'
$ws.Range('E14').Value = '  +8.06%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '685.40'
$ws.Range('E15').Value = '  +1.02%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '9.08'
$ws.Range('E16').Value = '  +5.21%  '
$ws.Range('D17').Value = '3.701.58'
$ws.Range('D18').Value = '71.770.53'
$ws.Range('E18').Value = '  +3.53%  '
$ws.Range('E19').Value = '  +2.31%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '18.22'
$ws.Range('E20').Value = '  +2.67%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.71'
$ws.Range('E21').Value = '  +3.47%  '
$ws.Range('B22').Value = 'Toncoin'
$ws.Range('C22').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.39'
$ws.Range('E22').Value = '  +19.01%  '
$ws.Range('B23').Value = 'Polygon'
$ws.Range('C23').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.947'
$ws.Range('E23').Value = '  +3.99%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '17.90'
$ws.Range('E24').Value = '  +5.23%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '103.98'
$ws.Range('E25').Value = '  +3.15%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '4.05'
$ws.Range('E26').Value = '  +4.10%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.86'
$ws.Range('E27').Value = '  +6.19%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '10.25'
$ws.Range('E28').Value = '  +6.38%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '35.43'
$ws.Range('E29').Value = '  +5.79%  '
$ws.Range('E30').Value = '  +6.34%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.40'
$ws.Range('E31').Value = '  +7.85%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.31'
$ws.Range('E32').Value = '  +16.41%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '11.33'
$ws.Range('E33').Value = '  +3.11%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '566.80'
$ws.Range('E34').Value = '  +2.97%  '
$ws.Range('E35').Value = '  +4.61%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '59.56'
$ws.Range('E36').Value = '  +2.64%  '
$ws.Range('D37').Value = '3.760.60'
$ws.Range('E37').Value = '  +4.15%  '
$ws.Range('E38').Value = '  -0.09%  '
$ws.Range('E39').Value = '  +3.82%  '
$ws.Range('D40').Value = '0.0₃0779'
$ws.Range('E40').Value = '  +6.14%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '35.71'
$ws.Range('E41').Value = '  +1.99%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0468'
$ws.Range('E43').Value = '  +10.29%  '
$ws.Range('E44').Value = '  +5.02%  '
$ws.Range('E45').Value = '  +6.08%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.91'
$ws.Range('E46').Value = '  +9.71%  '
$ws.Range('E47').Value = '  -0.38%  '
$ws.Range('E48').Value = '  +4.26%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.45'
$ws.Range('E49').Value = '  +3.37%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.999'
$ws.Range('E50').Value = '  -0.25%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '135.50'
